$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2024-05-21 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-22 Wednesday", 2)

# Update the division-problem table cells (addressed by row/column so that
# values which coincide with other cells' old/new text are not mismatched)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "676÷5=135, 1"
$t.Cell(1, 2).Range.Text = "110÷9=12, 2"
$t.Cell(1, 3).Range.Text = "319÷7=45, 4"
$t.Cell(1, 4).Range.Text = "102÷7=14, 4"
$t.Cell(1, 5).Range.Text = "223÷8=27, 7"

$t.Cell(5, 1).Range.Text = "893÷6=148, 5"
$t.Cell(5, 2).Range.Text = "786÷8=98, 2"
$t.Cell(5, 3).Range.Text = "649÷2=324, 1"
$t.Cell(5, 4).Range.Text = "414÷2=207, 0"
$t.Cell(5, 5).Range.Text = "490÷8=61, 2"

$t.Cell(9, 1).Range.Text = "951÷6=158, 3"
$t.Cell(9, 2).Range.Text = "778÷8=97, 2"
$t.Cell(9, 3).Range.Text = "584÷6=97, 2"
$t.Cell(9, 4).Range.Text = "547÷8=68, 3"
$t.Cell(9, 5).Range.Text = "314÷9=34, 8"

$t.Cell(13, 1).Range.Text = "977÷5=195, 2"
$t.Cell(13, 2).Range.Text = "713÷6=118, 5"
$t.Cell(13, 3).Range.Text = "155÷9=17, 2"
$t.Cell(13, 4).Range.Text = "741÷7=105, 6"
$t.Cell(13, 5).Range.Text = "378÷4=94, 2"

$t.Cell(17, 1).Range.Text = "717÷5=143, 2"
$t.Cell(17, 2).Range.Text = "342÷2=171, 0"
$t.Cell(17, 3).Range.Text = "625÷2=312, 1"
$t.Cell(17, 4).Range.Text = "586÷8=73, 2"
$t.Cell(17, 5).Range.Text = "251÷5=50, 1"

